$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    # Find.Execute's replacement argument runs AutoCorrect (e.g. smart quotes),
    # which would mangle literal apostrophes. Locate with Find (no replacement
    # argument) and then assign .Text on the collapsed/found range instead -
    # that path leaves the literal characters alone.
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText)
    if (-not $found) {
        throw "Text not found: $oldText"
    }
    $rng.Text = $newText
}

# 1. Header contact line: "Rochester, NY" -> "Open To Relocation & Remote Work"
Replace-ExactText "Rochester, NY | (727) 262-7305 | " "Open To Relocation & Remote Work | (727) 262-7305 | "

# 2. Summary paragraph: "with 1 year" -> "with over 1 year"
Replace-ExactText "Software Engineer with 1 year of experience" "Software Engineer with over 1 year of experience"

# 3. Job title: "Software Engineer (Contractor) " -> "Software Engineer "
Replace-ExactText "Software Engineer (Contractor) " "Software Engineer "

# 4. Job meta line: location + dates update
Replace-ExactText "| Lowe's Home Improvement (Revature) | Rochester, NY | Feb 2023 – Present (1 Year)" "| Lowe's Home Improvement (Revature) | Remote, US | Feb 2023 – May 2024"

# 5. Documentation bullet rewrite
Replace-ExactText "Created comprehensive documentation for onboarding new engineers, setting up development environments, and using internal tools." "Authored comprehensive documentation on development environment setup and internal tools, reducing new engineer onboarding time by 80% (from 10 days to 2 days)."

# 6. Remove "minor" before "bugs"
Replace-ExactText "by fixing minor bugs in a large codebase" "by fixing bugs in a large codebase"

# 7. Wireframes bullet rewrite
Replace-ExactText "Designed wireframes for a small web development team." "Created wireframes for web applications, collaborating with a cross-functional development team to understand design implementation feasibility."

# 8. Merge the two Honors & Awards scholarship bullets into one.
# Find the "Florida Bright Futures Scholarship" paragraph and delete the paragraph
# break right after it - that merges it into the next paragraph, adopting the next
# paragraph's pPr (which lacks the w:spacing element), matching the target shape.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "Florida Bright Futures Scholarship") {
        $p = $d.Paragraphs($i)
        $markRng = $d.Range($p.Range.End - 1, $p.Range.End)
        $markRng.Delete()
        break
    }
}
Replace-ExactText "Florida Bright Futures ScholarshipDoorways/Take Stock in Children Scholarship" "Florida Bright Futures and Doorways/Take Stock in Children Scholarships"

Write-Output "all edits applied"
